# Add the four new notes/questions captured after the farmer_rho_demo.bash
# entry, matching the new shared-string entries and sheet rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B50").Value = "grad_cost_and_rho vs. Find_Rho"
$ws.Range("D50").Value = "grad_cost_and_rho writes files; Find_Rho reads them"

$ws.Range("B52").Value = "grad_extension.mid_iter: are things computed before being written?"

$ws.Range("B54").Value = "The test py file needs better file setup and cleanup (e.g., get rid of _out files)"

# Match the author's final selection/cursor position after typing the new rows.
$ws.Range("B54").Select() | Out-Null
